$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sort the data range A1:G60 by column C (confidence) descending, using the
# worksheet Sort object so Excel persists a <sortState> in the saved XML
# (matching a Data > Sort operation performed through the UI).
$sortObj = $ws.Sort
$sortObj.SortFields.Clear()
$sortObj.SortFields.Add($ws.Range("C2:C60"), 0, 2, 0, 0)
$sortObj.SetRange($ws.Range("A1:G60"))
$sortObj.Header = 1
$sortObj.Apply()

# Widen column A to fit the longer account names now visible near the top.
$ws.Columns("A").ColumnWidth = 42.14

# Move/restore the active selection to A3.
$ws.Range("A3").Select() | Out-Null
